$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 'Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg'
$ws.Range("G3").Value2 = 'Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G4").Value2 = 'eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G5").Value2 = 'Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G6").Value2 = 'Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Range("G7").Value2 = 'NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg'
$ws.Range("G8").Value2 = 'AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg'
$ws.Range("G11").Value2 = 'Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Range("G12").Value2 = 'dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg'
$ws.Range("G13").Value2 = 'yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg'
$ws.Range("G17").Value2 = 'mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg'
$ws.Range("G24").Value2 = 'youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg'
$ws.Range("G25").Value2 = 'Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg'
$ws.Range("G27").Value2 = 'nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg'
$ws.Range("G28").Value2 = 'maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'
$ws.Range("G30").Value2 = 'wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
